$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the coin price (D) and 1h change (E) columns with the latest
# scraped figures. A few rows also moved in the ranking: BabyDogeCoin
# dropped out of the top listing and USDD entered, shifting Cronos,
# EnergySwap and Mantle up one row each, so the name/link cells (B/C)
# for rows 48-51 are refreshed as well.
#
# Some "D" prices (e.g. 212.77) are plain decimals that Excel would
# otherwise auto-convert to a floating point number (losing the exact
# source text to binary rounding). A leading apostrophe is used for
# those so they are stored as text, matching the original workbook.

$ws.Range("D2").Value = '26.469.65'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '1.626.74'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = "'212.77"
$ws.Range("E5").Value = '  -0.14%  '
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("D10").Value = "'18.78"
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = '1.853.56'
$ws.Range("D13").Value = '1.613.20'
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").Value = "'4.12"
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").Value = "'64.77"
$ws.Range("E16").Value = '  +2.68%  '
$ws.Range("D17").Value = '26.516.15'
$ws.Range("E17").Value = '  -0.56%  '
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = "'214.22"
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("E24").Value = '  +2.84%  '
$ws.Range("D25").Value = "'148.61"
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("D29").Value = "'15.51"
$ws.Range("E29").Value = '  +0.72%  '
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("E32").Value = '  +2.74%  '
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  -0.59%  '
$ws.Range("D35").Value = '1.216.98'
$ws.Range("E35").Value = '  +3.93%  '
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("E37").Value = '  +3.71%  '
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("E39").Value = '  -1.86%  '
$ws.Range("D40").Value = "'0.505"
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("E41").Value = '  -2.62%  '
$ws.Range("D42").Value = "'0.791"
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("D43").Value = "'5.35"
$ws.Range("E43").Value = '  -0.40%  '
$ws.Range("D44").Value = '1.764.94'
$ws.Range("E44").Value = '  -0.44%  '
$ws.Range("D45").Value = "'92.82"
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").Value = "'54.79"
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = "'0.0509"
$ws.Range("E48").Value = '  -0.53%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'7.52"
$ws.Range("E49").Value = '  -0.25%  '
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = '  +0.28%  '
